# "changes in concise marksheet" — update Corr/total marks on the quiz sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = 5
$ws.Range("B12").Value = 125
$ws.Range("E12").Value = "125/140"
